$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" (column F) values between row 16 and row 29
$val16 = $ws.Range("F16").Value2
$val29 = $ws.Range("F29").Value2

$ws.Range("F16").Value2 = $val29
$ws.Range("F29").Value2 = $val16
